# Remove the "Special Query" slide from the template.
#
# In the source presentation this is the last slide (slide 9 / sldId 308),
# which also owns notes page notesSlide7.xml. Deleting the slide through
# the Slides collection removes the slide part and cascades to remove its
# associated notes page as well.
$p = $ppt.ActivePresentation

for ($i = $p.Slides.Count; $i -ge 1; $i--) {
    $slide = $p.Slides.Item($i)
    $isSpecialQuery = $false

    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTextFrame) {
            if ($shape.TextFrame.TextRange.Text -eq "Special Query") {
                $isSpecialQuery = $true
                break
            }
        }
    }

    if ($isSpecialQuery) {
        $slide.Delete()
    }
}
